$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 17 (2023-02-27... date 44978): Oscar (G) 1 -> 3; Ines (H) -> 2; David (I) -> 2; Daniel (J) -> 2
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 2

# Row 19 (date 44980): Oscar (G) -> 2; Ines (H) -> 2; David (I) -> 2; Daniel (J) -> 2
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = 2

# Row 21 (date 44982): Oscar (G) -> 2
$ws.Range("G21").Value = 2

# Update the active cell selection to match the final state (I29)
$ws.Range("I29").Select()

$wb.Save()
